$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.503.84'
$ws.Range("E2").Value = '  +0.54%  '
$ws.Range("D3").Value = '2.108.00'
$ws.Range("E3").Value = '  +4.83%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '329.58'
$ws.Range("E5").Value = '  +1.49%  '
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5261'
$ws.Range("E7").Value = '  +2.57%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4366'
$ws.Range("E8").Value = '  +2.36%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08870'
$ws.Range("E9").Value = '  +2.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '47.16'
$ws.Range("E10").Value = '  +9.18%  '
$ws.Range("E11").Value = '  +2.37%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.53'
$ws.Range("E12").Value = '  -0.78%  '
$ws.Range("D13").Value = '2.108.87'
$ws.Range("E13").Value = '  +4.87%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.736'
$ws.Range("E14").Value = '  +2.61%  '
$ws.Range("E15").Value = '  +3.94%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '96.45'
$ws.Range("E16").Value = '  +2.23%  '
$ws.Range("E17").Value = '  -0.04%  '
$ws.Range("E18").Value = '  +1.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06642'
$ws.Range("E19").Value = '  +1.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.03'
$ws.Range("E20").Value = '  +0.80%  '
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("E22").Value = '  +2.22%  '
$ws.Range("D23").Value = '30.543.45'
$ws.Range("E23").Value = '  +0.49%  '
$ws.Range("E24").Value = '  +4.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.338'
$ws.Range("E25").Value = '  +3.45%  '
$ws.Range("D26").Value = '2.352.13'
$ws.Range("E26").Value = '  +4.74%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.43'
$ws.Range("E27").Value = '  -0.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.589'
$ws.Range("E28").Value = '  +6.79%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '162.03'
$ws.Range("E29").Value = '  -0.15%  '
$ws.Range("E30").Value = '  +1.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.213'
$ws.Range("E31").Value = '  +6.76%  '
$ws.Range("E32").Value = '  +2.36%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.675'
$ws.Range("E33").Value = '  +22.35%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.207'
$ws.Range("E34").Value = '  +2.31%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.923'
$ws.Range("E35").Value = '  +2.51%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.01'
$ws.Range("E36").Value = '  +9.89%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02588'
$ws.Range("E37").Value = '  +2.64%  '
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.486'
$ws.Range("E38").Value = '  +0.26%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06698'
$ws.Range("E39").Value = '  +0.30%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.71'
$ws.Range("E40").Value = '  +2.99%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2283'
$ws.Range("E41").Value = '  +4.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6802'
$ws.Range("E42").Value = '  +2.53%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.274'
$ws.Range("E43").Value = '  +2.78%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.04'
$ws.Range("E45").Value = '  +2.65%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6389'
$ws.Range("E46").Value = '  +3.59%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.212'
$ws.Range("E47").Value = '  +1.58%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.628'
$ws.Range("E48").Value = '  -0.87%  '
$ws.Range("E49").Value = '  -0.62%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '82.70'
$ws.Range("E50").Value = '  +2.63%  '
$ws.Range("B51").Value = 'WEMIXTOKEN'
$ws.Range("C51").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.196'
$ws.Range("E51").Value = '  +8.30%  '
